# Auto-generated edit script: updates Leve profit-calculation values
# across the per-class Leve sheets (ALC/BSM/CRP/CUL/GSM/LTW/WVR),
# reflecting refreshed Universalis market-price data.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 330.9
$ws.Range("I28").Value = 389.33334
$ws.Range("J28").Value = 243.25
$ws.Range("K28").Value = 389.33334
$ws.Range("L28").Value = 243.25
$ws.Range("M28").Value = 95.66665999999998
$ws.Range("N28").Value = -1213.25
$ws.Range("H86").Value = 3984.6667
$ws.Range("J86").Value = 4749
$ws.Range("L86").Value = 4749
$ws.Range("N86").Value = -6995
$ws.Range("H89").Value = 3984.6667
$ws.Range("J89").Value = 4749
$ws.Range("L89").Value = 23745
$ws.Range("N89").Value = -34977
$ws.Range("H106").Value = 23178.363
$ws.Range("I106").Value = 29495.25
$ws.Range("J106").Value = 6333.3335
$ws.Range("K106").Value = 29495.25
$ws.Range("L106").Value = 6333.3335
$ws.Range("M106").Value = -28864.25
$ws.Range("N106").Value = -7595.3335
$ws.Range("H138").Value = 4093.5967
$ws.Range("I138").Value = 1567.7778
$ws.Range("J138").Value = 4522.5093
$ws.Range("K138").Value = 4703.3334
$ws.Range("L138").Value = 13567.5279
$ws.Range("M138").Value = 436.6665999999996
$ws.Range("N138").Value = -23847.5279
$ws.Range("H141").Value = 3756.818
$ws.Range("I141").Value = 3619.4443
$ws.Range("K141").Value = 10858.3329
$ws.Range("M141").Value = -5678.332900000001

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3041.8
$ws.Range("I86").Value = 1301.5
$ws.Range("K86").Value = 1301.5
$ws.Range("M86").Value = -178.5
$ws.Range("H89").Value = 3041.8
$ws.Range("I89").Value = 1301.5
$ws.Range("K89").Value = 6507.5
$ws.Range("M89").Value = -891.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 571.8095
$ws.Range("H31").Value = 3956
$ws.Range("I31").Value = 2210.077
$ws.Range("K31").Value = 2210.077
$ws.Range("M31").Value = -1915.077
$ws.Range("H34").Value = 3956
$ws.Range("I34").Value = 2210.077
$ws.Range("K34").Value = 2210.077
$ws.Range("M34").Value = -2008.077
$ws.Range("H58").Value = 3214.3684
$ws.Range("I58").Value = 1773.7273
$ws.Range("K58").Value = 1773.7273
$ws.Range("M58").Value = -1570.7273
$ws.Range("H99").Value = 14786.15
$ws.Range("I99").Value = 10106.889
$ws.Range("K99").Value = 10106.889
$ws.Range("M99").Value = -8608.888999999999
$ws.Range("H113").Value = 571.8095
$ws.Range("H126").Value = 14786.15
$ws.Range("I126").Value = 10106.889
$ws.Range("K126").Value = 30320.667
$ws.Range("M126").Value = -27850.667
$ws.Range("H132").Value = 2500
$ws.Range("I132").Value = 1909.4231
$ws.Range("K132").Value = 5728.2693
$ws.Range("M132").Value = -3198.2693
$ws.Range("H136").Value = 3214.3684
$ws.Range("I136").Value = 1773.7273
$ws.Range("K136").Value = 5321.1819
$ws.Range("M136").Value = -2771.1819

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 44
$ws.Range("I11").Value = 51.75
$ws.Range("K11").Value = 155.25
$ws.Range("M11").Value = -15.25
$ws.Range("H69").Value = 3503
$ws.Range("J69").Value = 3000
$ws.Range("L69").Value = 9000
$ws.Range("N69").Value = -10622
$ws.Range("H72").Value = 3503
$ws.Range("J72").Value = 3000
$ws.Range("L72").Value = 27000
$ws.Range("N72").Value = -35112
$ws.Range("H105").Value = 7916.6665
$ws.Range("J105").Value = 7916.6665
$ws.Range("L105").Value = 23749.9995
$ws.Range("N105").Value = -28991.9995
$ws.Range("H111").Value = 2333.3333
$ws.Range("I111").Value = 1500
$ws.Range("J111").Value = 4000
$ws.Range("K111").Value = 4500
$ws.Range("L111").Value = 12000
$ws.Range("M111").Value = -1433
$ws.Range("N111").Value = -18134
$ws.Range("H120").Value = 11030
$ws.Range("J120").Value = 11479.311
$ws.Range("L120").Value = 34437.933
$ws.Range("N120").Value = -44113.933
$ws.Range("H122").Value = 821.55554
$ws.Range("I122").Value = 674.375
$ws.Range("K122").Value = 6069.375
$ws.Range("M122").Value = -3619.375
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("H139").Value = 2792.9375
$ws.Range("I139").Value = 2792.9375
$ws.Range("K139").Value = 8378.8125
$ws.Range("M139").Value = -3238.8125
$ws.Range("M130").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3139.0625
$ws.Range("J113").Value = 4999.75
$ws.Range("L113").Value = 4999.75
$ws.Range("N113").Value = -9339.75
$ws.Range("H122").Value = 1104974.8
$ws.Range("I122").Value = 205152.6
$ws.Range("K122").Value = 615457.8
$ws.Range("M122").Value = -613007.8
$ws.Range("H126").Value = 4974.75
$ws.Range("I126").Value = 4949
$ws.Range("J126").Value = 4983.3335
$ws.Range("K126").Value = 14847
$ws.Range("L126").Value = 14950.0005
$ws.Range("M126").Value = -12377
$ws.Range("N126").Value = -19890.0005

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1299
$ws.Range("I22").Value = 1299
$ws.Range("K22").Value = 1299
$ws.Range("M22").Value = -1004
$ws.Range("H27").Value = 1299
$ws.Range("I27").Value = 1299
$ws.Range("K27").Value = 1299
$ws.Range("M27").Value = -1192
$ws.Range("H46").Value = 3085.7144
$ws.Range("I46").Value = 1650
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 1650
$ws.Range("L46").Value = 5000
$ws.Range("M46").Value = -1462
$ws.Range("N46").Value = -5376
$ws.Range("H122").Value = 8001
$ws.Range("I122").Value = 10002
$ws.Range("J122").Value = 3999
$ws.Range("K122").Value = 30006
$ws.Range("L122").Value = 11997
$ws.Range("M122").Value = -27556
$ws.Range("N122").Value = -16897

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 12875
$ws.Range("I4").Value = 50000
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 50000
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -49887
$ws.Range("N4").Value = -726
$ws.Range("H100").Value = 1900.091
$ws.Range("I100").Value = 1655.6666
$ws.Range("K100").Value = 3311.3332
$ws.Range("M100").Value = -2770.3332
$ws.Range("H136").Value = 57420.89
$ws.Range("I136").Value = 1278.7693
$ws.Range("K136").Value = 3836.3079
$ws.Range("M136").Value = -1286.3079

